$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the stray row 13 (it held only the professor's name in B/C, with no
# label in column A); everything below shifts up one row, carrying its
# formatting (row heights, styles) with it.
$ws.Rows("13").Delete()

# --- Fix up the cell values that ended up different from a plain shift ---

# Name: / Environmental management in company — normalize the non-breaking
# space to a regular space.
$ws.Range("B4").Value = "Environmental management in company"
$ws.Range("C4").Value = "Environmental management in company"

# Objetivos: now shows the responsible professor instead of the long
# objectives paragraph.
$ws.Range("B10").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C10").Value = "5817650 - Érica Leonor Romão"

# Programa resumido: now shows "Semestral" instead of the short syllabus text.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: now shows the activation date instead of the long program text.
# Copy from B8/C8 (which already hold this exact text) so it stays a plain
# shared string instead of being reinterpreted as a date value.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

# Método: now shows the responsible professor instead of the evaluation
# method text.
$ws.Range("B18").Value = "5817650 - Érica Leonor Romão"
$ws.Range("C18").Value = "5817650 - Érica Leonor Romão"

# Critério:, Norma de recuperação: and Bibliografia: keep the text that used
# to belong to Método:, Critério: and Norma de recuperação: respectively
# (i.e. they did not shift together with the labels).
$ws.Range("B19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."
$ws.Range("C19").Value = "As avaliações serão por meio de trabalhos em equipes ou provas individuais, conforme adequação ao conteúdo programático."

$ws.Range("B20").Value = "Média ponderada das notas atribuídas às provas, exercícios, trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas às provas, exercícios, trabalhos práticos e relatórios."

$ws.Range("B21").Value = "A nota final será composta pela média obtida da nota do período somada à nota de recuperação"
$ws.Range("C21").Value = "A nota final será composta pela média obtida da nota do período somada à nota de recuperação"
